$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert the new "attribute" rows ---
# One new row right after Entity1 (row 3), before Entity2 (row 4).
$ws.Rows.Item(4).Insert()

# Three new rows right after Entity3 (after the inserts above, Entity3 now
# sits at row 9), before the total row (now row 10).
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# --- New row 4: Attr1 under Entity1 ---
$ws.Range("A6").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Attr1"

# --- New rows 10-12: Attr1/Attr2/Attr3 under Entity3 ---
$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Attr1"

$ws.Range("A7").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Attr2"

$ws.Range("A8").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Attr3"

$excel.CutCopyMode = 0

# Row 12 (Attr3 under Entity3) has no C cell at all, unlike its siblings.
$ws.Range("C12").Clear()

# Reflect the final cursor position, like the saved workbook shows.
$ws.Range("A14").Select()
